$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Append the new day's row (2025-11-03) right after the last existing row (28).
# Force the date-like label to be stored as literal text (matching the export's
# "Date" column, which stores every date as a string, not a real Excel date).
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "2025-11-03"
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 115
